# Apply cell updates from the crypto price refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells hold price strings that must stay as exact text
# (dotted thousands separators, trailing zeros, small-decimal values).
# Force the cell format to Text before assigning so Excel does not
# re-interpret the string as a number and lose precision/formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.915.12"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.279.75"
$ws.Range("E3").Value = "  -4.13%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.87"
$ws.Range("E5").Value = "  -4.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.11"
$ws.Range("E6").Value = "  -7.97%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.281.58"
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.464"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.118"
$ws.Range("E11").Value = "  -5.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.405"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.842.33"
$ws.Range("E13").Value = "  -4.14%  "
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.84"
$ws.Range("E15").Value = "  -6.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.273.34"
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("E17").Value = "  -4.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.059.72"
$ws.Range("E18").Value = "  -3.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -6.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.75"
$ws.Range("E20").Value = "  -5.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.53"
$ws.Range("E21").Value = "  -4.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "373.16"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.68"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.532"
$ws.Range("E25").Value = "  -6.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.423.92"
$ws.Range("E26").Value = "  -3.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000101"
$ws.Range("E27").Value = "  -10.16%  "
$ws.Range("E28").Value = "  -5.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.10"
$ws.Range("E30").Value = "  -7.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -5.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.46"
$ws.Range("E33").Value = "  -5.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.50"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.24"
$ws.Range("E35").Value = "  -7.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.05"
$ws.Range("E36").Value = "  -7.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.26"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.51"
$ws.Range("E38").Value = "  -6.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.62"
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.309.30"
$ws.Range("E40").Value = "  -4.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.98"
$ws.Range("E41").Value = "  -16.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0725"
$ws.Range("E42").Value = "  -7.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.66"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.745"
$ws.Range("E44").Value = "  -4.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.09"
$ws.Range("E45").Value = "  -7.03%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.56"
$ws.Range("E46").Value = "  -6.87%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.10"
$ws.Range("E47").Value = "  -6.31%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.328.53"
$ws.Range("E49").Value = "  -8.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.35"
$ws.Range("E50").Value = "  -7.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.20"
$ws.Range("E51").Value = "  -6.11%  "
